$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename sheets
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Sheet1").Name = "not sure about this"
$wb.Worksheets.Item("error codes").Name = "Error Code XML"

$ws1 = $wb.Worksheets.Item("not sure about this")
$ws2 = $wb.Worksheets.Item("Error Code XML")

# ---------------------------------------------------------------------------
# 2) Rebuild the "Error Code XML" sheet content
# ---------------------------------------------------------------------------

# Drop the old summary row (row 18) - it's no longer part of the sheet.
$ws2.Rows.Item(18).Delete()

# New header row
$ws2.Range("A1").Value = "Error Scenario"
$ws2.Range("B1").Value = "XML File Content"
$ws2.Range("A1:B1").HorizontalAlignment = -4108   # xlCenter

# Row 2 - invalid user credentials
$ws2.Range("A2").Value = "Scraper attempts to scrape with invalid user credentials"
$ws2.Range("B2").Value = "<scrape-session>`n<base-url>www.elen7045.co.za</base-url>`n<date>01/06/2009</>`n<time>12:59:34</time>`n<error>InvalidCredentials</error>`n"
$ws2.Rows.Item(2).RowHeight = 74.25
$ws2.Range("A2").VerticalAlignment = -4108        # xlCenter
$ws2.Range("B2").VerticalAlignment = -4108
$ws2.Range("B2").WrapText = $true

# Row 3 - customer hasn't signed up for e-billing
$ws2.Range("A3").Value = "Scraper attempts to scrape with customer that hasn" + [char]8217 + "t signed up for e-billing"
$ws2.Range("B3").Value = "<scrape-session>`n<base-url>www.elen7045.co.za</base-url>`n<date>01/06/2009</>`n<time>12:59:34</time>`n<error>UserNotSignedUpForEbilling</error>`n"
$ws2.Rows.Item(3).RowHeight = 79.5
$ws2.Range("A3:B3").VerticalAlignment = -4108
$ws2.Range("A3:B3").WrapText = $true

# Row 4 - e-billing process required further customer input
$ws2.Range("A4").Value = "Scraper attempts to scrape with customer who" + [char]8217 + "s e-billing process required further customer input to complete"
$ws2.Range("B4").Value = "<scrape-session>`n<base-url>www.elen7045.co.za</base-url>`n<date>01/06/2009</>`n<time>12:59:34</time>`n<error>AccountUpdateRequired</error>`n"
$ws2.Rows.Item(4).RowHeight = 75
$ws2.Range("A4:B4").VerticalAlignment = -4108
$ws2.Range("A4:B4").WrapText = $true

# Row 5 - billing company website is down
$ws2.Range("A5").Value = "Scraper attempts to scrape while the billing company website is down"
$ws2.Range("B5").Value = "<scrape-session>`n<base-url>www.elen7045.co.za</base-url>`n<date>01/06/2009</>`n<time>12:59:34</time>`n<error>BillingSiteDown</error>`n"
$ws2.Rows.Item(5).RowHeight = 75.75
$ws2.Range("A5:B5").VerticalAlignment = -4108
$ws2.Range("A5:B5").WrapText = $true

# Row 6 - e-billing service is unavailable
$ws2.Range("A6").Value = "Scraper attempts to scrape while e-billing service is unavailable"
$ws2.Range("B6").Value = "<scrape-session>`n<base-url>www.elen7045.co.za</base-url>`n<date>01/06/2009</>`n<time>12:59:34</time>`n<error>BillingSitePageError</error>`n"
$ws2.Rows.Item(6).RowHeight = 77.25
$ws2.Range("A6").WrapText = $true
$ws2.Range("B6").VerticalAlignment = -4108
$ws2.Range("B6").WrapText = $true

# Row 7 - nonconforming script
$ws2.Range("A7").Value = "Scraper attempts to scrape with a nonconforming script"
$ws2.Range("B7").Value = "<scrape-session>`n<base-url>www.elen7045.co.za</base-url>`n<date>01/06/2009</>`n<time>12:59:34</time>`n<error>UndeterminedError</error>`n"
$ws2.Rows.Item(7).RowHeight = 80.25
$ws2.Range("A7:B7").VerticalAlignment = -4108
$ws2.Range("A7:B7").WrapText = $true

# Column widths for the new layout
$ws2.Columns.Item(1).ColumnWidth = 50.5
$ws2.Columns.Item(2).ColumnWidth = 41

# ---------------------------------------------------------------------------
# 3) View / selection bookkeeping to mirror the authored state
# ---------------------------------------------------------------------------
$ws1.Range("B18").Select()
$ws2.Range("A7").Select()

$ws1.Select()
